$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header in B1 from "Cabang (Khusus Non-BCA)*" to "Cabang*"
$ws.Range("B1").Value = "Cabang*"

# Update the selected cell to B2 (matches the post-edit cursor position)
$ws.Range("B2").Select()
